# [FIX] Nueva actualizacion del driverchrome
#
# The "FLAG" column (T) on the main "Sheet" tab was bulk-updated from "SI"
# to "NO" for every data row (rows 2-65; the later rows already contained
# other flag values / were left untouched upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
$ws.Activate()

$ws.Range("T2:T65").Value = "NO"

# Leave the selection where the author ended up before saving.
$ws.Range("AA67").Select()
